$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change yes/no answers to "Yes" (reuses the existing "Yes" shared string)
$ws.Range("B13").Value = "Yes"
$ws.Range("B27").Value = "Yes"
$ws.Range("B33").Value = "Yes"
$ws.Range("B43").Value = "Yes"
$ws.Range("B84").Value = "Yes"
$ws.Range("B87").Value = "Yes"
$ws.Range("B88").Value = "Yes"
$ws.Range("B98").Value = "Yes"

# Fill in the relevant line-number / test-name references in column E, in the
# order the new unique strings were originally authored so the shared-string
# table ends up laid out the same way.
$ws.Range("E84").Value = "Line 64"
$ws.Range("E83").Value = "Line 63"
$ws.Range("E80").Value = "Lines 57-60"
$ws.Range("E33").Value = "DetectionTests/threeSuitTest"
$ws.Range("E27").Value = "DetectionTests/oneOffRoyalTest"
$ws.Range("E98").Value = "Line 67"

# These reuse existing shared strings ("Line 61" / "Line 62")
$ws.Range("E81").Value = "Line 61"
$ws.Range("E82").Value = "Line 62"

# Plain numeric values (no shared string involved)
$ws.Range("E87").Value = 66
$ws.Range("E88").Value = 65

# Update the view state (scroll position / active selection)
$ws.Application.ActiveWindow.ScrollRow = 101
$ws.Range("B85").Select()
